$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D and E columns to remain text so numeric-looking values
# (e.g. "9.66") are not auto-converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '70.297.70'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '3.604.91'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '578.58'
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").Value = '189.52'
$ws.Range("E6").Value = '  -3.19%  '
$ws.Range("D7").Value = '3.602.39'
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("D8").Value = '0.629'
$ws.Range("E8").Value = '  -3.05%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").Value = '0.189'
$ws.Range("E10").Value = '  +4.00%  '
$ws.Range("D11").Value = '0.657'
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("D12").Value = '55.97'
$ws.Range("E12").Value = '  -4.74%  '
$ws.Range("D13").Value = '0.0000314'
$ws.Range("E13").Value = '  +6.24%  '
$ws.Range("D14").Value = '9.66'
$ws.Range("E14").Value = '  -3.35%  '
$ws.Range("D15").Value = '4.186.98'
$ws.Range("E15").Value = '  -0.96%  '
$ws.Range("D16").Value = '19.82'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").Value = '3.605.20'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = '70.356.18'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = '12.68'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("E21").Value = '  -2.99%  '
$ws.Range("D22").Value = '493.59'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").Value = '19.42'
$ws.Range("E23").Value = '  +2.51%  '
$ws.Range("D24").Value = '4.92'
$ws.Range("E24").Value = '  -8.33%  '
$ws.Range("D25").Value = '96.03'
$ws.Range("E25").Value = '  +4.53%  '
$ws.Range("D26").Value = '4.36'
$ws.Range("E26").Value = '  -3.28%  '
$ws.Range("D27").Value = '2.99'
$ws.Range("E27").Value = '  -5.88%  '
$ws.Range("D28").Value = '11.14'
$ws.Range("E28").Value = '  -3.91%  '
$ws.Range("D29").Value = '9.40'
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("D30").Value = '32.13'
$ws.Range("E30").Value = '  -3.19%  '
$ws.Range("D31").Value = '7.64'
$ws.Range("E31").Value = '  -4.00%  '
$ws.Range("D32").Value = '12.17'
$ws.Range("E32").Value = '  -1.38%  '
$ws.Range("D33").Value = '66.26'
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("E34").Value = '  -3.94%  '
$ws.Range("D35").Value = '574.30'
$ws.Range("D36").Value = '38.57'
$ws.Range("E36").Value = '  -5.42%  '
$ws.Range("D37").Value = '0.0₃0812'
$ws.Range("E37").Value = '  -4.45%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '3.06'
$ws.Range("E39").Value = '  +5.53%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '0.396'
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '3.32'
$ws.Range("E41").Value = '  +11.99%  '
$ws.Range("D42").Value = '3.56'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  -7.05%  '
$ws.Range("D44").Value = '3.02'
$ws.Range("E44").Value = '  -5.64%  '
$ws.Range("D45").Value = '3.235.54'
$ws.Range("E45").Value = '  -2.67%  '
$ws.Range("D46").Value = '9.84'
$ws.Range("E46").Value = '  +5.86%  '
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("D48").Value = '3.42'
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '3.19'
$ws.Range("E51").Value = '  -4.45%  '

# Restore the default cell style (removes the temporary text
# number-format override so no stray style index is left behind).
$ws.Range("D2:E51").Style = "Normal"

